$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 44595
$ws.Range("L2").Value = "Primera"
$ws.Range("N2").Value = 15500
$ws.Range("O2").Value = 16000
$ws.Range("P2").Value = 15750
$ws.Range("S2").Value = 875
$ws.Range("D3").Value = 44634
$ws.Range("L3").Value = "Especial"
$ws.Range("M3").Value = 200
$ws.Range("N3").Value = 14000
$ws.Range("O3").Value = 15000
$ws.Range("P3").Value = 14500
$ws.Range("S3").Value = 806
$ws.Range("D4").Value = 44634
$ws.Range("L4").Value = "Primera"
$ws.Range("M4").Value = 200
$ws.Range("N4").Value = 10000
$ws.Range("P4").Value = 10500
$ws.Range("S4").Value = 583
$ws.Range("D5").Value = 44630
$ws.Range("L5").Value = "Especial"
$ws.Range("M5").Value = 300
$ws.Range("N5").Value = 15000
$ws.Range("O5").Value = 16000
$ws.Range("P5").Value = 15500
$ws.Range("S5").Value = 861
$ws.Range("D6").Value = 44630
$ws.Range("L6").Value = "Primera"
$ws.Range("M6").Value = 300
$ws.Range("N6").Value = 12000
$ws.Range("O6").Value = 13000
$ws.Range("P6").Value = 12500
$ws.Range("S6").Value = 694
$ws.Range("D7").Value = 44630
$ws.Range("L7").Value = "Segunda"
$ws.Range("N7").Value = 9000
$ws.Range("O7").Value = 10000
$ws.Range("P7").Value = 9500
$ws.Range("S7").Value = 528
$ws.Range("D8").Value = 44603
$ws.Range("L8").Value = "Especial"
$ws.Range("M8").Value = 240
$ws.Range("N8").Value = 14500
$ws.Range("O8").Value = 15000
$ws.Range("P8").Value = 14750
$ws.Range("S8").Value = 819
$ws.Range("D9").Value = 44636
$ws.Range("L9").Value = "Especial"
$ws.Range("M9").Value = 240
$ws.Range("N9").Value = 14000
$ws.Range("O9").Value = 15000
$ws.Range("P9").Value = 14500
$ws.Range("S9").Value = 806
$ws.Range("D10").Value = 44636
$ws.Range("L10").Value = "Primera"
$ws.Range("N10").Value = 10000
$ws.Range("O10").Value = 11000
$ws.Range("P10").Value = 10500
$ws.Range("S10").Value = 583
$ws.Range("D11").Value = 44295
$ws.Range("L11").Value = "Especial"
$ws.Range("M11").Value = 200
$ws.Range("N11").Value = 14500
$ws.Range("O11").Value = 15000
$ws.Range("P11").Value = 14750
$ws.Range("S11").Value = 819
$ws.Range("D12").Value = 44295
$ws.Range("L12").Value = "Primera"
$ws.Range("M12").Value = 200
$ws.Range("N12").Value = 12500
$ws.Range("O12").Value = 13000
$ws.Range("P12").Value = 12750
$ws.Range("S12").Value = 708
$ws.Range("D13").Value = 44295
$ws.Range("L13").Value = "Segunda"
$ws.Range("M13").Value = 240
$ws.Range("N13").Value = 10500
$ws.Range("P13").Value = 10750
$ws.Range("S13").Value = 597
$ws.Range("D16").Value = 44609
$ws.Range("L16").Value = "Primera"
$ws.Range("M16").Value = 240
$ws.Range("N16").Value = 13000
$ws.Range("O16").Value = 14000
$ws.Range("P16").Value = 13500
$ws.Range("S16").Value = 750
$ws.Range("D17").Value = 44609
$ws.Range("L17").Value = "Segunda"
$ws.Range("M17").Value = 240
$ws.Range("N17").Value = 11000
$ws.Range("O17").Value = 12000
$ws.Range("P17").Value = 11500
$ws.Range("S17").Value = 639
$ws.Range("D18").Value = 44631
$ws.Range("L18").Value = "Especial"
$ws.Range("N18").Value = 15000
$ws.Range("O18").Value = 16000
$ws.Range("P18").Value = 15500
$ws.Range("S18").Value = 861
$ws.Range("D19").Value = 44631
$ws.Range("M19").Value = 248
$ws.Range("N19").Value = 12000
$ws.Range("O19").Value = 13000
$ws.Range("P19").Value = 12516
$ws.Range("S19").Value = 695
$ws.Range("D20").Value = 44631
$ws.Range("M20").Value = 200
$ws.Range("D21").Value = 44637
$ws.Range("N21").Value = 14000
$ws.Range("P21").Value = 14500
$ws.Range("S21").Value = 806
$ws.Range("D22").Value = 44637
$ws.Range("M22").Value = 240
$ws.Range("N22").Value = 10000
$ws.Range("O22").Value = 11000
$ws.Range("P22").Value = 10500
$ws.Range("S22").Value = 583
$ws.Range("D23").Value = 44294
$ws.Range("L23").Value = "Especial"
$ws.Range("M23").Value = 200
$ws.Range("N23").Value = 14500
$ws.Range("O23").Value = 15000
$ws.Range("P23").Value = 14750
$ws.Range("S23").Value = 819
$ws.Range("D24").Value = 44294
$ws.Range("M24").Value = 240
$ws.Range("N24").Value = 12500
$ws.Range("O24").Value = 13000
$ws.Range("P24").Value = 12750
$ws.Range("S24").Value = 708
$ws.Range("D25").Value = 44294
$ws.Range("L25").Value = "Segunda"
$ws.Range("M25").Value = 240
$ws.Range("N25").Value = 10500
$ws.Range("O25").Value = 11000
$ws.Range("P25").Value = 10750
$ws.Range("S25").Value = 597
$ws.Range("D26").Value = 44606
$ws.Range("M26").Value = 240
$ws.Range("N26").Value = 11500
$ws.Range("O26").Value = 12000
$ws.Range("P26").Value = 11750
$ws.Range("S26").Value = 653
$ws.Range("D27").Value = 44606
$ws.Range("L27").Value = "Segunda"
$ws.Range("N27").Value = 9500
$ws.Range("O27").Value = 10000
$ws.Range("P27").Value = 9750
$ws.Range("S27").Value = 542
$ws.Range("D28").Value = 44607
$ws.Range("M28").Value = 300
$ws.Range("N28").Value = 11000
$ws.Range("O28").Value = 12000
$ws.Range("P28").Value = 11500
$ws.Range("S28").Value = 639
$ws.Range("D29").Value = 44607
$ws.Range("N29").Value = 9000
$ws.Range("O29").Value = 10000
$ws.Range("P29").Value = 9500
$ws.Range("S29").Value = 528
